$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (2025-10-03) was still blank (style "Normal with border" on every
# cell). Fill it in the same way rows 2 and 3 already were: most cells get
# the "Good" cell style, the "Personnalisation" column (B) and one of the
# "Castes / Éduc." entries (F) get the "Neutral" style, keeping the
# existing thin border + centered alignment that the named Excel cell
# styles alone don't carry.
#
# Row 2 already has exactly this style layout (B=Neutral, C=Good, D=Bad,
# E..J=Good, K..M=Bad) except F, which needs to be Neutral instead of
# Good, so: copy row 2's formatting onto row 4, then patch F4's format
# from B2 (a cell that is already styled "Neutral").
$xlPasteFormats = -4122

$ws.Range("B2:M2").Copy() | Out-Null
$ws.Range("B4:M4").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("B2").Copy() | Out-Null
$ws.Range("F4").PasteSpecial($xlPasteFormats) | Out-Null

$excel.CutCopyMode = $false

# Now fill in the actual values (entered in the same order the author
# originally typed them).
$ws.Range("B4").Value = "Criminal"
$ws.Range("C4").Value = "Charge"
$ws.Range("G4").Value = "Artisanat"
$ws.Range("H4").Value = "Cuirassé"
$ws.Range("J4").Value = "Élémentariste"
$ws.Range("I4").Value = "Armes de finesse"
$ws.Range("F4").Value = "Apprenti d’un maître"
$ws.Range("E4").Value = "Aigrefin"
